# Added results for partial model
# Updates the sliding-window predicted values (IPC PO), DELTA, and
# DELTA^2 columns (plus a few float-precision-only IPC RO updates)
# for the window-10 partial-model results sheet, plus the TOTAL/MSE
# summary rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 30.01235008239746
$ws.Cells.Item(2,4).Value = 0.09235008239745923
$ws.Cells.Item(2,5).Value = 0.00852853771881751
$ws.Cells.Item(3,3).Value = 29.95737266540527
$ws.Cells.Item(3,4).Value = -0.02262733459473054
$ws.Cells.Item(3,5).Value = 0.0005119962708618896
$ws.Cells.Item(4,2).Value = 30.03999999999999
$ws.Cells.Item(4,3).Value = 30.14963340759277
$ws.Cells.Item(4,4).Value = 0.1096334075927814
$ws.Cells.Item(4,5).Value = 0.01201948406040494
$ws.Cells.Item(5,2).Value = 30.21000000000001
$ws.Cells.Item(5,3).Value = 30.0929012298584
$ws.Cells.Item(5,4).Value = -0.1170987701416095
$ws.Cells.Item(5,5).Value = 0.0137121219686775
$ws.Cells.Item(6,3).Value = 30.26950645446777
$ws.Cells.Item(6,4).Value = 0.04950645446777457
$ws.Cells.Item(6,5).Value = 0.002450889033969837
$ws.Cells.Item(7,3).Value = 30.25530052185059
$ws.Cells.Item(7,4).Value = -0.1246994781494095
$ws.Cells.Item(7,5).Value = 0.01554995985073506
$ws.Cells.Item(8,3).Value = 30.50997734069824
$ws.Cells.Item(8,4).Value = 0.06997734069824446
$ws.Cells.Item(8,5).Value = 0.00489682821119818
$ws.Cells.Item(9,3).Value = 30.38161277770996
$ws.Cells.Item(9,4).Value = -0.09838722229004304
$ws.Cells.Item(9,5).Value = 0.009680045509950342
$ws.Cells.Item(10,3).Value = 30.44911766052246
$ws.Cells.Item(10,4).Value = -0.2408823394775368
$ws.Cells.Item(10,5).Value = 0.05802430147217128
$ws.Cells.Item(11,3).Value = 30.4470043182373
$ws.Cells.Item(11,4).Value = -0.3029956817626953
$ws.Cells.Item(11,5).Value = 0.09180638316684053
$ws.Cells.Item(12,3).Value = 30.63541984558105
$ws.Cells.Item(12,4).Value = -0.304580154418943
$ws.Cells.Item(12,5).Value = 0.09276907046586719
$ws.Cells.Item(13,3).Value = 30.7780818939209
$ws.Cells.Item(13,4).Value = -0.1719181060791044
$ws.Cells.Item(13,5).Value = 0.0295558351978262
$ws.Cells.Item(14,3).Value = 31.21818351745605
$ws.Cells.Item(14,4).Value = 0.1981835174560587
$ws.Cells.Item(14,5).Value = 0.03927670659125591
$ws.Cells.Item(15,3).Value = 31.37577629089355
$ws.Cells.Item(15,4).Value = 0.2557762908935501
$ws.Cells.Item(15,5).Value = 0.06542151098326197
$ws.Cells.Item(16,3).Value = 31.4897289276123
$ws.Cells.Item(16,4).Value = 0.2097289276123036
$ws.Cells.Item(16,5).Value = 0.04398622307740686
$ws.Cells.Item(17,3).Value = 31.24315452575684
$ws.Cells.Item(17,4).Value = -0.1368454742431595
$ws.Cells.Item(17,5).Value = 0.01872668382083523
$ws.Cells.Item(18,3).Value = 31.55834007263184
$ws.Cells.Item(18,4).Value = -0.02165992736816236
$ws.Cells.Item(18,5).Value = 0.0004691524535940687
$ws.Cells.Item(19,2).Value = 31.65000000000001
$ws.Cells.Item(19,3).Value = 31.96619987487793
$ws.Cells.Item(19,4).Value = 0.316199874877924
$ws.Cells.Item(19,5).Value = 0.0999823608728148
$ws.Cells.Item(20,3).Value = 32.57534408569336
$ws.Cells.Item(20,4).Value = 0.6953440856933639
$ws.Cells.Item(20,5).Value = 0.4835033975087402
$ws.Cells.Item(21,3).Value = 32.38410568237305
$ws.Cells.Item(21,4).Value = 0.1041056823730457
$ws.Cells.Item(21,5).Value = 0.01083799310235749
$ws.Cells.Item(22,3).Value = 32.5141487121582
$ws.Cells.Item(22,4).Value = 0.06414871215820028
$ws.Cells.Item(22,5).Value = 0.004115057271555633
$ws.Cells.Item(23,2).Value = 32.84999999999999
$ws.Cells.Item(23,3).Value = 32.70607376098633
$ws.Cells.Item(23,4).Value = -0.1439262390136662
$ws.Cells.Item(23,5).Value = 0.02071476227661897
$ws.Cells.Item(24,2).Value = 32.90000000000001
$ws.Cells.Item(24,3).Value = 32.94353866577148
$ws.Cells.Item(24,4).Value = 0.04353866577147869
$ws.Cells.Item(24,5).Value = 0.00189561541716053
$ws.Cells.Item(25,2).Value = 33.09999999999999
$ws.Cells.Item(25,3).Value = 32.95751571655273
$ws.Cells.Item(25,4).Value = -0.1424842834472599
$ws.Cells.Item(25,5).Value = 0.02030177102947911
$ws.Cells.Item(26,2).Value = 33.40000000000001
$ws.Cells.Item(26,3).Value = 33.68034744262695
$ws.Cells.Item(26,4).Value = 0.2803474426269474
$ws.Cells.Item(26,5).Value = 0.07859468858746958
$ws.Cells.Item(27,3).Value = 33.67167282104492
$ws.Cells.Item(27,4).Value = -0.02832717895508097
$ws.Cells.Item(27,5).Value = 0.000802429067553182
$ws.Cells.Item(28,2).Value = 34.09999999999999
$ws.Cells.Item(28,3).Value = 33.86410903930664
$ws.Cells.Item(28,4).Value = -0.2358909606933537
$ws.Cells.Item(28,5).Value = 0.05564454533683334
$ws.Cells.Item(29,2).Value = 34.40000000000001
$ws.Cells.Item(29,3).Value = 34.45993423461914
$ws.Cells.Item(29,4).Value = 0.05993423461913494
$ws.Cells.Item(29,5).Value = 0.003592112479381513
$ws.Cells.Item(30,2).Value = 34.90000000000001
$ws.Cells.Item(30,3).Value = 35.07357406616211
$ws.Cells.Item(30,4).Value = 0.1735740661621037
$ws.Cells.Item(30,5).Value = 0.03012795644404635
$ws.Cells.Item(31,3).Value = 35.75087356567383
$ws.Cells.Item(31,4).Value = 0.450873565673831
$ws.Cells.Item(31,5).Value = 0.2032869722234344
$ws.Cells.Item(32,3).Value = 36.01235580444336
$ws.Cells.Item(32,4).Value = 0.3123558044433565
$ws.Cells.Item(32,5).Value = 0.09756614856945639
$ws.Cells.Item(33,3).Value = 36.00337600708008
$ws.Cells.Item(33,4).Value = -0.296623992919919
$ws.Cells.Item(33,5).Value = 0.08798579317575618
$ws.Cells.Item(34,3).Value = 36.55116653442383
$ws.Cells.Item(34,4).Value = -0.248833465576169
$ws.Cells.Item(34,5).Value = 0.0619180935906465
$ws.Cells.Item(35,3).Value = 37.06188583374023
$ws.Cells.Item(35,4).Value = -0.2381141662597628
$ws.Cells.Item(35,5).Value = 0.05669835617358195
$ws.Cells.Item(36,2).Value = 37.90000000000001
$ws.Cells.Item(36,3).Value = 37.8508415222168
$ws.Cells.Item(36,4).Value = -0.04915847778320881
$ws.Cells.Item(36,5).Value = 0.002416555937962234
$ws.Cells.Item(37,3).Value = 38.3694953918457
$ws.Cells.Item(37,4).Value = -0.1305046081542969
$ws.Cells.Item(37,5).Value = 0.01703145274950657
$ws.Cells.Item(38,2).Value = 38.90000000000001
$ws.Cells.Item(38,3).Value = 39.00495529174805
$ws.Cells.Item(38,4).Value = 0.1049552917480412
$ws.Cells.Item(38,5).Value = 0.01101561326591644
$ws.Cells.Item(39,2).Value = 39.40000000000001
$ws.Cells.Item(39,3).Value = 39.52904891967773
$ws.Cells.Item(39,4).Value = 0.1290489196777287
$ws.Cells.Item(39,5).Value = 0.01665362366998887
$ws.Cells.Item(40,2).Value = 39.90000000000001
$ws.Cells.Item(40,3).Value = 39.70820999145508
$ws.Cells.Item(40,4).Value = -0.1917900085449276
$ws.Cells.Item(40,5).Value = 0.03678340737766338
$ws.Cells.Item(41,2).Value = 40.09999999999999
$ws.Cells.Item(41,3).Value = 39.92831802368164
$ws.Cells.Item(41,4).Value = -0.1716819763183537
$ws.Cells.Item(41,5).Value = 0.02947470099257576
$ws.Cells.Item(42,2).Value = 40.59999999999999
$ws.Cells.Item(42,3).Value = 40.23454666137695
$ws.Cells.Item(42,4).Value = -0.3654533386230412
$ws.Cells.Item(42,5).Value = 0.1335561427107272
$ws.Cells.Item(43,2).Value = 40.90000000000001
$ws.Cells.Item(43,3).Value = 40.47024154663086
$ws.Cells.Item(43,4).Value = -0.4297584533691463
$ws.Cells.Item(43,5).Value = 0.1846923282422407
$ws.Cells.Item(44,2).Value = 41.20000000000001
$ws.Cells.Item(44,3).Value = 41.21378707885742
$ws.Cells.Item(44,4).Value = 0.01378707885741193
$ws.Cells.Item(44,5).Value = 0.000190083543420495
$ws.Cells.Item(45,3).Value = 41.2248649597168
$ws.Cells.Item(45,4).Value = -0.2751350402832031
$ws.Cells.Item(45,5).Value = 0.07569929039163981
$ws.Cells.Item(46,3).Value = 41.653076171875
$ws.Cells.Item(46,4).Value = -0.1469238281249972
$ws.Cells.Item(46,5).Value = 0.0215866112709037
$ws.Cells.Item(47,3).Value = 42.15799331665039
$ws.Cells.Item(47,4).Value = -0.04200668334961222
$ws.Cells.Item(47,5).Value = 0.001764561446034588
$ws.Cells.Item(48,3).Value = 43.88340759277344
$ws.Cells.Item(48,4).Value = 1.183407592773435
$ws.Cells.Item(48,5).Value = 1.400453530633815
$ws.Cells.Item(49,3).Value = 44.43407821655273
$ws.Cells.Item(49,4).Value = 0.7340782165527386
$ws.Cells.Item(49,5).Value = 0.5388708280172494
$ws.Cells.Item(50,3).Value = 44.31099700927734
$ws.Cells.Item(50,4).Value = 0.1109970092773409
$ws.Cells.Item(50,5).Value = 0.0123203360685141
$ws.Cells.Item(51,3).Value = 44.46113967895508
$ws.Cells.Item(51,4).Value = -1.138860321044923
$ws.Cells.Item(51,5).Value = 1.297002830850546
$ws.Cells.Item(52,3).Value = -0.05531524658206166
$ws.Cells.Item(52,5).Value = 5.604465680149265
$ws.Cells.Item(53,5).Value = 0.1120893136029853

Write-Output "Updated sliding window results (window 10) for partial model."
